$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(0.353672031788087, 0.2465263520751981, 0.3418284350654049, 0.3112793583715512, 0.7152945399284363, 0.7077127695083618, 0.9207399487495422, 0.8079603910446167)

for ($row = 2; $row -le 26; $row++) {
    for ($col = 2; $col -le 9; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 2]
    }
}
